# Auto-generated PowerShell COM-interop script applying the KHL injuries update
# for 2025-10-27: refresh snapshot timestamps, move "Стюарт Дин" (СЮЛ) to the
# "returned" sheet, and add newly-injured "Уткин Дмитрий А" (СОЧ) to snapshot +
# the "new_injured" sheet.

$wb = $excel.ActiveWorkbook

# --- Sheet "snapshot": refresh date/timestamp columns, and update the player
#     rows for teams СОЧ/СПР/СЮЛ (rows 30-38) to reflect the roster shift caused
#     by Стюарт Дин returning and Уткин Дмитрий А becoming newly injured.
$wsSnapshot = $wb.Worksheets.Item("snapshot")

$wsSnapshot.Range("A2").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L2").Value = "2025-10-27T20:28:22.871346"

$wsSnapshot.Range("A3").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L3").Value = "2025-10-27T20:28:26.453559"

$wsSnapshot.Range("A4").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L4").Value = "2025-10-27T20:28:26.453586"

$wsSnapshot.Range("A5").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L5").Value = "2025-10-27T20:28:26.453604"

$wsSnapshot.Range("A6").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L6").Value = "2025-10-27T20:28:26.453621"

$wsSnapshot.Range("A7").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L7").Value = "2025-10-27T20:28:26.453637"

$wsSnapshot.Range("A8").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L8").Value = "2025-10-27T20:28:35.815926"

$wsSnapshot.Range("A9").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L9").Value = "2025-10-27T20:28:35.815955"

$wsSnapshot.Range("A10").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L10").Value = "2025-10-27T20:28:38.716488"

$wsSnapshot.Range("A11").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L11").Value = "2025-10-27T20:28:41.253774"

$wsSnapshot.Range("A12").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L12").Value = "2025-10-27T20:28:41.253803"

$wsSnapshot.Range("A13").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L13").Value = "2025-10-27T20:28:43.845087"

$wsSnapshot.Range("A14").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L14").Value = "2025-10-27T20:28:43.845114"

$wsSnapshot.Range("A15").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L15").Value = "2025-10-27T20:28:43.845132"

$wsSnapshot.Range("A16").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L16").Value = "2025-10-27T20:28:43.845150"

$wsSnapshot.Range("A17").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L17").Value = "2025-10-27T20:28:52.142197"

$wsSnapshot.Range("A18").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L18").Value = "2025-10-27T20:28:55.082295"

$wsSnapshot.Range("A19").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L19").Value = "2025-10-27T20:28:57.625787"

$wsSnapshot.Range("A20").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L20").Value = "2025-10-27T20:29:00.125266"

$wsSnapshot.Range("A21").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L21").Value = "2025-10-27T20:29:00.125294"

$wsSnapshot.Range("A22").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L22").Value = "2025-10-27T20:29:03.018397"

$wsSnapshot.Range("A23").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L23").Value = "2025-10-27T20:29:03.018423"

$wsSnapshot.Range("A24").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L24").Value = "2025-10-27T20:29:03.018439"

$wsSnapshot.Range("A25").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L25").Value = "2025-10-27T20:29:03.018456"

$wsSnapshot.Range("A26").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L26").Value = "2025-10-27T20:29:03.018472"

$wsSnapshot.Range("A27").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L27").Value = "2025-10-27T20:29:05.483522"

$wsSnapshot.Range("A28").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L28").Value = "2025-10-27T20:29:10.743538"

$wsSnapshot.Range("A29").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L29").Value = "2025-10-27T20:29:10.743562"

$wsSnapshot.Range("A30").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E30").Value = "Уткин Дмитрий А"
$wsSnapshot.Range("F30").Value = "15"
$wsSnapshot.Range("G30").Value = "нападающий"
$wsSnapshot.Range("H30").Value = "35195"
$wsSnapshot.Range("I30").Value = "1369_СОЧ_уткиндмитрийа"
$wsSnapshot.Range("L30").Value = "2025-10-27T20:29:10.743579"

$wsSnapshot.Range("A31").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("B31").Value = "СОЧ"
$wsSnapshot.Range("C31").Value = "ХК Сочи"
$wsSnapshot.Range("D31").Value = "hc_sochi"
$wsSnapshot.Range("E31").Value = "Хомченко Павел"
$wsSnapshot.Range("F31").Value = "30"
$wsSnapshot.Range("G31").Value = "вратарь"
$wsSnapshot.Range("H31").Value = "17592"
$wsSnapshot.Range("I31").Value = "1369_СОЧ_хомченкопавел"
$wsSnapshot.Range("K31").Value = "https://www.khl.ru/clubs/hc_sochi/team/"
$wsSnapshot.Range("L31").Value = "2025-10-27T20:29:10.743595"

$wsSnapshot.Range("A32").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("B32").Value = "СПР"
$wsSnapshot.Range("C32").Value = "Спартак"
$wsSnapshot.Range("D32").Value = "spartak"
$wsSnapshot.Range("E32").Value = "Порядин Павел"
$wsSnapshot.Range("F32").Value = "24"
$wsSnapshot.Range("H32").Value = "19258"
$wsSnapshot.Range("I32").Value = "1369_СПР_порядинпавел"
$wsSnapshot.Range("K32").Value = "https://www.khl.ru/clubs/spartak/team/"
$wsSnapshot.Range("L32").Value = "2025-10-27T20:29:13.263160"

$wsSnapshot.Range("A33").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E33").Value = "Алалыкин Данил"
$wsSnapshot.Range("F33").Value = "61"
$wsSnapshot.Range("H33").Value = "34493"
$wsSnapshot.Range("I33").Value = "1369_СЮЛ_алалыкинданил"
$wsSnapshot.Range("L33").Value = "2025-10-27T20:29:15.799036"

$wsSnapshot.Range("A34").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E34").Value = "Берлёв Антон"
$wsSnapshot.Range("F34").Value = "83"
$wsSnapshot.Range("G34").Value = "нападающий"
$wsSnapshot.Range("H34").Value = "20546"
$wsSnapshot.Range("I34").Value = "1369_СЮЛ_берлевантон"
$wsSnapshot.Range("L34").Value = "2025-10-27T20:29:15.799060"

$wsSnapshot.Range("A35").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E35").Value = "Зоркин Никита"
$wsSnapshot.Range("F35").Value = "52"
$wsSnapshot.Range("G35").Value = "защитник"
$wsSnapshot.Range("H35").Value = "26738"
$wsSnapshot.Range("I35").Value = "1369_СЮЛ_зоркинникита"
$wsSnapshot.Range("L35").Value = "2025-10-27T20:29:15.799076"

$wsSnapshot.Range("A36").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E36").Value = "Кузьмин Глеб"
$wsSnapshot.Range("F36").Value = "17"
$wsSnapshot.Range("H36").Value = "22170"
$wsSnapshot.Range("I36").Value = "1369_СЮЛ_кузьминглеб"
$wsSnapshot.Range("L36").Value = "2025-10-27T20:29:15.799092"

$wsSnapshot.Range("A37").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E37").Value = "Набиев Артём"
$wsSnapshot.Range("F37").Value = "79"
$wsSnapshot.Range("H37").Value = "41187"
$wsSnapshot.Range("I37").Value = "1369_СЮЛ_набиевартем"
$wsSnapshot.Range("L37").Value = "2025-10-27T20:29:15.799107"

$wsSnapshot.Range("A38").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("E38").Value = "Пименов Артём"
$wsSnapshot.Range("F38").Value = "68"
$wsSnapshot.Range("G38").Value = "нападающий"
$wsSnapshot.Range("H38").Value = "21205"
$wsSnapshot.Range("I38").Value = "1369_СЮЛ_пименовартем"
$wsSnapshot.Range("L38").Value = "2025-10-27T20:29:15.799124"

$wsSnapshot.Range("A39").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L39").Value = "2025-10-27T20:29:15.799140"

$wsSnapshot.Range("A40").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L40").Value = "2025-10-27T20:29:15.799156"

$wsSnapshot.Range("A41").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L41").Value = "2025-10-27T20:29:15.799169"

$wsSnapshot.Range("A42").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L42").Value = "2025-10-27T20:29:18.255366"

$wsSnapshot.Range("A43").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L43").Value = "2025-10-27T20:29:18.255391"

$wsSnapshot.Range("A44").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L44").Value = "2025-10-27T20:29:23.772661"

$wsSnapshot.Range("A45").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L45").Value = "2025-10-27T20:29:26.244275"

$wsSnapshot.Range("A46").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L46").Value = "2025-10-27T20:29:26.244300"

$wsSnapshot.Range("A47").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L47").Value = "2025-10-27T20:29:26.244317"

$wsSnapshot.Range("A48").Value = "по состоянию на 27 октября 2025"
$wsSnapshot.Range("L48").Value = "2025-10-27T20:29:26.244332"

# --- Sheet "returned": log Стюарт Дин (СЮЛ) as returned from injury.
$wsReturned = $wb.Worksheets.Item("returned")
$wsReturned.Range("A2").Value = "СЮЛ"
$wsReturned.Range("B2").Value = "Салават Юлаев"
$wsReturned.Range("C2").Value = "Стюарт Дин"
$wsReturned.Range("D2").Value = "1369_СЮЛ_стюартдин"
$wsReturned.Range("E2").Value = "RETURN"
$wsReturned.Range("F2").Value = "2025-10-27T20:29:26.849817"

# --- Sheet "new_injured": log Уткин Дмитрий А (СОЧ) as newly injured.
$wsNewInjured = $wb.Worksheets.Item("new_injured")
$wsNewInjured.Range("A2").Value = "СОЧ"
$wsNewInjured.Range("B2").Value = "ХК Сочи"
$wsNewInjured.Range("C2").Value = "Уткин Дмитрий А"
$wsNewInjured.Range("D2").Value = "1369_СОЧ_уткиндмитрийа"
$wsNewInjured.Range("E2").Value = "INJURED_NEW"
$wsNewInjured.Range("F2").Value = "2025-10-27T20:29:26.853860"

